# Professional_Competence_Model.pptx edit
# - resize/reposition the three roundRect boxes ("Knowledge", "Professional vision",
#   "Observable behavior" -> "Student-directed behavior")
# - update the text of the third box
# - reposition/resize the two connectors between the boxes (one now flipped vertically)

function ToPt($emu) {
    # Convert EMU -> points for the Shape.Left/Top/Width/Height properties.
    # A small epsilon compensates for floating point truncation inside the
    # host so that the EMU value written back to the XML matches exactly.
    return ($emu / 12700.0) + 0.00003
}

function Set-ShapeXfrm($shape, $x, $y, $cx, $cy) {
    $shape.Left = ToPt($x)
    $shape.Top = ToPt($y)
    $shape.Width = ToPt($cx)
    $shape.Height = ToPt($cy)
}

function Get-Chars($tr, $a, $b) {
    # Select the 0-based half-open character range [a, b) of $tr as a
    # Characters() sub-range, working around the host's off-by-one handling
    # of the (start, length) arguments of TextRange.Characters().
    if ($a -eq 0) {
        return $tr.Characters(1, $b + 1)
    } else {
        return $tr.Characters($a + 2, $b - $a)
    }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=4 "Textfeld 3" (Knowledge) ---
$shKnowledge = $s.Shapes.Item(1)
Set-ShapeXfrm $shKnowledge 132178 1775220 3224629 2145268

# --- Shape id=5 "Textfeld 4" (Professional vision) ---
$shVision = $s.Shapes.Item(2)
Set-ShapeXfrm $shVision 3910428 1777758 4331872 2145268

# --- Shape id=6 "Textfeld 5" (Observable behavior -> Student-directed behavior) ---
$shBehavior = $s.Shapes.Item(3)

$tr = $shBehavior.TextFrame.TextRange
# The heading paragraph currently reads "Observable behavior". The word
# "Observable" (characters [0,10)) becomes "Student-directed", keeping the
# following " behavior" untouched. A second pass re-writes the "directed"
# sub-range so that it ends up as its own run (matching the run split in
# the target markup). This is done before resizing the shape because the
# textbox has auto-fit enabled, which would otherwise override an
# explicitly-set height as soon as the text reflows.
(Get-Chars $tr 0 10).Text = "Student-directed"
(Get-Chars $tr 8 16).Text = "directed"

Set-ShapeXfrm $shBehavior 8795921 1735054 3224629 2213372

# --- Connector id=3 (Knowledge -> Professional vision) ---
$cxn1 = $s.Shapes.Item(4)
Set-ShapeXfrm $cxn1 3356807 2847854 553621 2538

# --- Connector id=10 (Professional vision -> Student-directed behavior) ---
$cxn2 = $s.Shapes.Item(5)
Set-ShapeXfrm $cxn2 8242300 2841740 553621 8652
$cxn2.VerticalFlip = -1
